$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Region"
$ws.Range("B1").Value = "Country"

$ws.Name = "region list"

$ws.Range("B2").Select() | Out-Null
